# Requirements Traceability Matrix (RTM) update:
# Remove the Defect ID references from the "Comments" column (J) on the
# RTM sheet, per the Release Agent's request to strip Defect IDs before
# uploading the RTM for Release 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTM")

# Rows in column J that currently hold a Defect ID (e.g. "DE458") or a
# "was impacted by DE### (discovered in TC####)" note - clear them all.
$rows = @(3,4,5,21,23,28,29,30,31,32,33,34,35,36,40)

foreach ($r in $rows) {
    $ws.Range("J$r").Value = $null
}
